$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values stay as text, matching the source data
# (values like "566.70", "2.00", "7.40" must not be coerced to numbers,
# which would silently drop the trailing zero).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.750.21'
$ws.Range("E2").Value = '  +2.85%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.436.79'
$ws.Range("E3").Value = '  +1.65%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '566.70'
$ws.Range("E5").Value = '  +2.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.25'
$ws.Range("E6").Value = '  +5.86%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("E8").Value = '  +2.41%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.172'
$ws.Range("E9").Value = '  +9.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '2.435.88'
$ws.Range("E10").Value = '  +1.75%  '
$ws.Range("E11").Value = '  -1.54%  '
$ws.Range("E12").Value = '  +2.66%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.71'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000179'
$ws.Range("E14").Value = '  +6.51%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '69.397.62'
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.879.58'
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.08'
$ws.Range("E17").Value = '  +5.86%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.438.54'
$ws.Range("E18").Value = '  +1.73%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.88'
$ws.Range("E19").Value = '  +6.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '344.36'
$ws.Range("E20").Value = '  +4.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.19'
$ws.Range("E21").Value = '  +5.99%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.88'
$ws.Range("E22").Value = '  +3.54%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.00'
$ws.Range("E23").Value = '  +7.79%  '
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '66.11'
$ws.Range("E25").Value = '  +0.83%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.83'
$ws.Range("E26").Value = '  +6.09%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.51'
$ws.Range("E27").Value = '  +5.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.556.97'
$ws.Range("E28").Value = '  +1.43%  '
$ws.Range("E29").Value = '  +0.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0857'
$ws.Range("E30").Value = '  +8.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.40'
$ws.Range("E31").Value = '  +5.32%  '
$ws.Range("E32").Value = '  +11.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '456.77'
$ws.Range("E33").Value = '  +10.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  +2.52%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.26'
$ws.Range("E36").Value = '  +1.34%  '
$ws.Range("E37").Value = '  +7.59%  '
$ws.Range("E38").Value = '  +0.71%  '
$ws.Range("E39").Value = '  +0.05%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.27'
$ws.Range("E40").Value = '  +3.65%  '
$ws.Range("E41").Value = '  +4.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.43'
$ws.Range("E42").Value = '  +5.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.53'
$ws.Range("E43").Value = '  +5.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '37.99'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("E45").Value = '  +3.57%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.11'
$ws.Range("E46").Value = '  +9.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '136.32'
$ws.Range("E47").Value = '  +6.35%  '
$ws.Range("E48").Value = '  +3.79%  '
$ws.Range("B49").Value = 'ARBITRUM'
$ws.Range("C49").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.492'
$ws.Range("E49").Value = '  +3.99%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0722'
$ws.Range("E50").Value = '  +2.49%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.565'
$ws.Range("E51").Value = '  +2.49%  '
